$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 319.8158
$ws.Range("I6").Value = 75.8
$ws.Range("K6").Value = 227.4
$ws.Range("M6").Value = -115.4
$ws.Range("H38").Value = 1980.2858
$ws.Range("I38").Value = 340.5
$ws.Range("J38").Value = 4166.6665
$ws.Range("K38").Value = 1021.5
$ws.Range("L38").Value = 12499.9995
$ws.Range("M38").Value = -649.5
$ws.Range("N38").Value = -13243.9995
$ws.Range("H40").Value = 4526.9165
$ws.Range("I40").Value = 3580.75
$ws.Range("K40").Value = 3580.75
$ws.Range("M40").Value = -3405.75
$ws.Range("H70").Value = 7454.4
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 7772.5713
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 23317.7139
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -23857.7139
$ws.Range("H73").Value = 7454.4
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 7772.5713
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 23317.7139
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -25189.7139
$ws.Range("H76").Value = 4996.25
$ws.Range("I76").Value = 5001.5
$ws.Range("K76").Value = 5001.5
$ws.Range("M76").Value = -4686.5
$ws.Range("H79").Value = 4996.25
$ws.Range("I79").Value = 5001.5
$ws.Range("K79").Value = 5001.5
$ws.Range("M79").Value = -3909.5
$ws.Range("H113").Value = 2966.4443
$ws.Range("J113").Value = 3959.8
$ws.Range("L113").Value = 3959.8
$ws.Range("N113").Value = -10467.8

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6888
$ws.Range("I32").Value = 6271.057
$ws.Range("K32").Value = 6271.057
$ws.Range("M32").Value = -5984.057
$ws.Range("H45").Value = 2072.7
$ws.Range("I45").Value = 2047.4375
$ws.Range("J45").Value = 2173.75
$ws.Range("K45").Value = 2047.4375
$ws.Range("L45").Value = 2173.75
$ws.Range("M45").Value = -1670.4375
$ws.Range("N45").Value = -2927.75
$ws.Range("H74").Value = 1958.3334
$ws.Range("I74").Value = 1409.1818
$ws.Range("K74").Value = 1409.1818
$ws.Range("M74").Value = -535.1818000000001
$ws.Range("H77").Value = 1958.3334
$ws.Range("I77").Value = 1409.1818
$ws.Range("K77").Value = 7045.909000000001
$ws.Range("M77").Value = -2677.909000000001
$ws.Range("H117").Value = 29999
$ws.Range("J117").Value = 29999
$ws.Range("L117").Value = 29999
$ws.Range("N117").Value = -39177

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2520.2354
$ws.Range("I20").Value = 2324.6667
$ws.Range("J20").Value = 2836.1538
$ws.Range("K20").Value = 2324.6667
$ws.Range("L20").Value = 2836.1538
$ws.Range("M20").Value = -2077.6667
$ws.Range("N20").Value = -3330.1538
$ws.Range("H80").Value = 640.375
$ws.Range("J80").Value = 726.2727
$ws.Range("L80").Value = 726.2727
$ws.Range("N80").Value = -2722.2727
$ws.Range("H83").Value = 640.375
$ws.Range("J83").Value = 726.2727
$ws.Range("L83").Value = 3631.3635
$ws.Range("N83").Value = -13615.3635
$ws.Range("H94").Value = 3345.923
$ws.Range("I94").Value = 2770.3333
$ws.Range("K94").Value = 2770.3333
$ws.Range("M94").Value = -2319.3333

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1468.3478
$ws.Range("J16").Value = 2278.4
$ws.Range("L16").Value = 2278.4
$ws.Range("N16").Value = -2852.4
$ws.Range("H31").Value = 3132.6155
$ws.Range("I31").Value = 1527.5625
$ws.Range("J31").Value = 5700.7
$ws.Range("K31").Value = 1527.5625
$ws.Range("L31").Value = 5700.7
$ws.Range("M31").Value = -1232.5625
$ws.Range("N31").Value = -6290.7
$ws.Range("H34").Value = 3132.6155
$ws.Range("I34").Value = 1527.5625
$ws.Range("J34").Value = 5700.7
$ws.Range("K34").Value = 1527.5625
$ws.Range("L34").Value = 5700.7
$ws.Range("M34").Value = -1325.5625
$ws.Range("N34").Value = -6104.7
$ws.Range("H70").Value = 43999
$ws.Range("J70").Value = 43999
$ws.Range("L70").Value = 43999
$ws.Range("N70").Value = -44629
$ws.Range("H73").Value = 43999
$ws.Range("J73").Value = 43999
$ws.Range("L73").Value = 43999
$ws.Range("N73").Value = -46183
$ws.Range("H113").Value = 1468.3478
$ws.Range("J113").Value = 2278.4
$ws.Range("L113").Value = 2278.4
$ws.Range("N113").Value = -6618.4
$ws.Range("H132").Value = 1602778.5
$ws.Range("I132").Value = 2002627.4
$ws.Range("K132").Value = 6007882.199999999
$ws.Range("M132").Value = -6005352.199999999
$ws.Range("H141").Value = 137794.83
$ws.Range("I141").Value = 40000
$ws.Range("J141").Value = 157353.8
$ws.Range("K141").Value = 40000
$ws.Range("L141").Value = 157353.8
$ws.Range("M141").Value = -34820
$ws.Range("N141").Value = -167713.8

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4183.0625
$ws.Range("I5").Value = 439.05884
$ws.Range("J5").Value = 8426.267
$ws.Range("K5").Value = 1317.17652
$ws.Range("L5").Value = 25278.801
$ws.Range("M5").Value = -1205.17652
$ws.Range("N5").Value = -25502.801
$ws.Range("H7").Value = 280.625
$ws.Range("I7").Value = 233
$ws.Range("J7").Value = 423.5
$ws.Range("K7").Value = 699
$ws.Range("L7").Value = 1270.5
$ws.Range("M7").Value = -587
$ws.Range("N7").Value = -1494.5
$ws.Range("H86").Value = 994.8889
$ws.Range("I86").Value = 988
$ws.Range("K86").Value = 2964
$ws.Range("M86").Value = -1778
$ws.Range("H89").Value = 994.8889
$ws.Range("I89").Value = 988
$ws.Range("K89").Value = 8892
$ws.Range("M89").Value = -2964
$ws.Range("H92").Value = 461.2857
$ws.Range("J92").Value = 585.8889
$ws.Range("L92").Value = 1757.6667
$ws.Range("N92").Value = -4253.6667
$ws.Range("H107").Value = 10581
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 10581
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 31743
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -35583
$ws.Range("H114").Value = 22729620
$ws.Range("J114").Value = 3172.625
$ws.Range("L114").Value = 9517.875
$ws.Range("N114").Value = -16025.875
$ws.Range("H131").Value = 1826.683
$ws.Range("J131").Value = 1868.3684
$ws.Range("L131").Value = 5605.1052
$ws.Range("N131").Value = -15685.1052
$ws.Range("H135").Value = 4183.0625
$ws.Range("I135").Value = 439.05884
$ws.Range("J135").Value = 8426.267
$ws.Range("K135").Value = 3951.52956
$ws.Range("L135").Value = 75836.40299999999
$ws.Range("M135").Value = -1416.52956
$ws.Range("N135").Value = -80906.40299999999

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4122.5625
$ws.Range("I80").Value = 3973.818
$ws.Range("K80").Value = 3973.818
$ws.Range("M80").Value = -2975.818
$ws.Range("H83").Value = 4122.5625
$ws.Range("I83").Value = 3973.818
$ws.Range("K83").Value = 19869.09
$ws.Range("M83").Value = -14877.09
$ws.Range("H122").Value = 1812.4166
$ws.Range("I122").Value = 1684.5
$ws.Range("J122").Value = 2452
$ws.Range("K122").Value = 5053.5
$ws.Range("L122").Value = 7356
$ws.Range("M122").Value = -2603.5
$ws.Range("N122").Value = -12256

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3116.1177
$ws.Range("I22").Value = 1648.3334
$ws.Range("J22").Value = 3430.6428
$ws.Range("K22").Value = 1648.3334
$ws.Range("L22").Value = 3430.6428
$ws.Range("M22").Value = -1353.3334
$ws.Range("N22").Value = -4020.6428
$ws.Range("H27").Value = 3116.1177
$ws.Range("I27").Value = 1648.3334
$ws.Range("J27").Value = 3430.6428
$ws.Range("K27").Value = 1648.3334
$ws.Range("L27").Value = 3430.6428
$ws.Range("M27").Value = -1541.3334
$ws.Range("N27").Value = -3644.6428
$ws.Range("H40").Value = 2528.28
$ws.Range("I40").Value = 2593.1904
$ws.Range("J40").Value = 2187.5
$ws.Range("K40").Value = 2593.1904
$ws.Range("L40").Value = 2187.5
$ws.Range("M40").Value = -2457.1904
$ws.Range("N40").Value = -2459.5
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4251
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 25000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -21256
$ws.Range("N71").ClearContents()
$ws.Range("H93").Value = 1564.3077
$ws.Range("I93").Value = 1564.3077
$ws.Range("K93").Value = 1564.3077
$ws.Range("M93").Value = -316.3077000000001
$ws.Range("H98").Value = 68500
$ws.Range("J98").Value = 68500
$ws.Range("L98").Value = 68500
$ws.Range("N98").Value = -74490
$ws.Range("H136").Value = 4171.387
$ws.Range("J136").Value = 5799.8335
$ws.Range("L136").Value = 17399.5005
$ws.Range("N136").Value = -22499.5005

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 89998
$ws.Range("J75").Value = 89998
$ws.Range("L75").Value = 89998
$ws.Range("N75").Value = -91870
$ws.Range("H78").Value = 89998
$ws.Range("J78").Value = 89998
$ws.Range("L78").Value = 269994
$ws.Range("N78").Value = -279354
$ws.Range("H125").Value = 99000
$ws.Range("J125").Value = 99000
$ws.Range("L125").Value = 99000
$ws.Range("N125").Value = -108840
$ws.Range("H135").Value = 82500
$ws.Range("J135").Value = 82500
$ws.Range("L135").Value = 82500
$ws.Range("N135").Value = -92640
$ws.Range("H136").Value = 24915.777
$ws.Range("I136").Value = 36407.168
$ws.Range("K136").Value = 109221.504
$ws.Range("M136").Value = -106671.504

Write-Output "Applied all Faerie_Profits updates"